$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows (for line7 / line8) before the current row 8 (extr1),
# shifting existing extr1..extr8 rows down from 8-15 to 10-17.
$ws.Range("A8:E9").Insert()

# --- New row 8: line7 ---
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# --- New row 9: line8 ---
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Match formatting of column A (bold, centered, bordered) used throughout the table
$ws.Range("A8:A9").Font.Bold = $true
$ws.Range("A8:A9").HorizontalAlignment = -4108
$ws.Range("A8:A9").VerticalAlignment = -4160
$ws.Range("A8:A9").Borders.LineStyle = 1

# --- Update shifted rows (previously rows 8-15, now rows 10-17) ---
# Column A must become the sequential index again (8..15)

# Row 10 (was row8 / extr1): A 6->8, C 10->5, D 11->12, E False->True
$ws.Range("A10").Value = 8
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# Row 11 (was row9 / extr2): A 7->9, C 7->5, D 8->9, E stays True
$ws.Range("A11").Value = 9
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# Row 12 (was row10 / extr3): A 8->10, C 9->10, D stays 11, E True->False
$ws.Range("A12").Value = 10
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $false

# Row 13 (was row11 / extr4): A 9->11, C stays 7, D 11->8, E True->False
$ws.Range("A13").Value = 11
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $false

# Row 14 (was row12 / extr5): A 10->12, C 5->9, D 7->11, E False->True
$ws.Range("A14").Value = 12
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $true

# Row 15 (was row13 / extr6): A 11->13, C 8->7, D new 11, E False->True
$ws.Range("A15").Value = 13
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $true

# Row 16 (was row14 / extr7): A 12->14, C 5, D 7, E True
$ws.Range("A16").Value = 14
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

# Row 17 (was row15 / extr8): A 13->15, C 8, D 5, E True
$ws.Range("A17").Value = 15
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
